$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 640
$ws1.Range("F3").Value = 6068
$ws1.Range("F6").Value = 1022
$ws1.Range("F7").Value = 408
$ws1.Range("F8").Value = 1385
$ws1.Range("F10").Value = 3137
$ws1.Range("F11").Value = 404
$ws1.Range("F12").Value = 1983
$ws1.Range("F13").Value = 124
$ws1.Range("F15").Value = 207
$ws1.Range("F16").Value = 90
$ws1.Range("F17").Value = 190
$ws1.Range("F18").Value = 1007
$ws1.Range("F19").Value = 371
$ws1.Range("F20").Value = 61
$ws1.Range("F21").Value = 93
$ws1.Range("F22").Value = 3748
$ws1.Range("F23").Value = 1187
$ws1.Range("F24").Value = 2978
$ws1.Range("F26").Value = 2532
$ws1.Range("F27").Value = 4320
$ws1.Range("F29").Value = 938
$ws1.Range("F31").Value = 1365
$ws1.Range("F32").Value = 136
$ws1.Range("F33").Value = 16
$ws1.Range("F34").Value = 48
$ws1.Range("F35").Value = 46
$ws1.Range("F37").Value = 1044
$ws1.Range("F38").Value = 1300
$ws1.Range("F40").Value = 1130
$ws1.Range("F41").Value = 732
$ws1.Range("F42").Value = 648
$ws1.Range("F43").Value = 441
$ws1.Range("F44").Value = 25
$ws1.Range("F45").Value = 129
$ws1.Range("F47").Value = 7
$ws1.Range("F49").Value = 3624

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 920
$ws2.Range("F25").Value = 21
$ws2.Range("F26").Value = 10

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 640
$ws4.Range("F3").Value = 6068
$ws4.Range("F6").Value = 408
$ws4.Range("F7").Value = 1385
$ws4.Range("F8").Value = 3137
$ws4.Range("F10").Value = 1983
$ws4.Range("F11").Value = 124
$ws4.Range("F13").Value = 207
$ws4.Range("F14").Value = 920
$ws4.Range("F16").Value = 90
$ws4.Range("F17").Value = 190
$ws4.Range("F18").Value = 1007
$ws4.Range("F19").Value = 371
$ws4.Range("F20").Value = 93
$ws4.Range("F21").Value = 3748
$ws4.Range("F23").Value = 1187
$ws4.Range("F25").Value = 2978
$ws4.Range("F26").Value = 2532
$ws4.Range("F27").Value = 4320
$ws4.Range("F30").Value = 938
$ws4.Range("F31").Value = 1365
$ws4.Range("F33").Value = 1044
$ws4.Range("F35").Value = 1300
$ws4.Range("F37").Value = 1130
$ws4.Range("F39").Value = 732
$ws4.Range("F41").Value = 441
$ws4.Range("F43").Value = 25
$ws4.Range("F44").Value = 21
$ws4.Range("F45").Value = 129
$ws4.Range("F46").Value = 10
$ws4.Range("F48").Value = 3624
